# Update attendance/visitor figures ("F" column) on the "展览" (Exhibitions)
# sheet and the merged "全部类型" (All types) sheet, per upstream data refresh
# (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value  = 13121
$wsExpo.Range("F4").Value  = 34
$wsExpo.Range("F5").Value  = 95
$wsExpo.Range("F6").Value  = 105
$wsExpo.Range("F8").Value  = 31
$wsExpo.Range("F10").Value = 13074
$wsExpo.Range("F11").Value = 312
$wsExpo.Range("F12").Value = 556
$wsExpo.Range("F13").Value = 8779
$wsExpo.Range("F14").Value = 7827
$wsExpo.Range("F16").Value = 134
$wsExpo.Range("F18").Value = 141
$wsExpo.Range("F19").Value = 996
$wsExpo.Range("F24").Value = 342
$wsExpo.Range("F26").Value = 5221

# --- Sheet "全部类型" (all types, merged view) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value  = 9
$wsAll.Range("F4").Value  = 13121
$wsAll.Range("F5").Value  = 34
$wsAll.Range("F6").Value  = 95
$wsAll.Range("F7").Value  = 105
$wsAll.Range("F9").Value  = 31
$wsAll.Range("F11").Value = 13074
$wsAll.Range("F12").Value = 312
$wsAll.Range("F13").Value = 556
$wsAll.Range("F14").Value = 8779
$wsAll.Range("F15").Value = 7827
$wsAll.Range("F17").Value = 134
$wsAll.Range("F19").Value = 141
$wsAll.Range("F20").Value = 996
$wsAll.Range("F23").Value = 9
$wsAll.Range("F27").Value = 342
$wsAll.Range("F29").Value = 5221
